$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.696.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.227.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.76'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.81%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.402'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0874'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.62%  '
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.559.97'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.795'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.233.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.582.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0887'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '246.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.62%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.138'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.58'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.76%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('E34').Value = '  +8.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0621'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('E37').Value = '  -4.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.000238'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +30.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.75%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.59%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0234'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0973'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.465.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.84%  '
$ws.Range('E51').Value = '  -0.85%  '
